# Moonscraper Chart Editor Manual.docx - apply commit:
#   "Fixed the Open Note FAQ point."
#
# The bulk of the underlying XML diff is Word re-flowing the document
# after its spell-checker pass (dropping stale <w:proofErr/> bookmarks,
# which merges adjacent same-formatted runs). The one real content
# change is the FAQ answer for "How do I place open notes?" - it used
# to say "press 6" and now explains the open-note key is "0" (with the
# 1-5 keys also toggling back to standard notes). Word's hidden
# "_GoBack" last-edit bookmark also moves to that edited sentence.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0. Relocate the hidden "_GoBack" bookmark: remove it from its old
#    (now stale) spot; it will be re-added after the FAQ text below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 1. Re-touch every run pair/triple that used to be split around a
#    <w:proofErr/> spell-check marker so Word folds them back into a
#    single run and drops the now-stale proofErr markers. None of
#    these change the visible text - replacement text equals the
#    search text.
# ------------------------------------------------------------------
function Clean-Text($text) {
    $d.Content.Find.Execute($text, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2) | Out-Null
}

Clean-Text "Moonscraper Chart Editor Manual"
Clean-Text "Developed by Alexander “FireFox” Ong"
Clean-Text "Welcome to Moonscraper Chart Editor"
Clean-Text "Extract all the contents of the rar file"
Clean-Text "If Moonscraper doesn’t respond"
Clean-Text "calibration of Moonscraper so that"
Clean-Text "Moonscraper divides chart editing"
Clean-Text "By default Moonscraper opens up to the local view"
Clean-Text "of a song, those events being notes and starpower."
Clean-Text "drag on a note to create a sustain"
Clean-Text "Starpower- Click to place a starpower event"
Clean-Text "restart Moonscraper with your controller"
Clean-Text "By default Moonscraper saves in the .chart format"
Clean-Text "However, Moonscraper can convert"
Clean-Text "You can change how Moonscraper saves the songs"
Clean-Text "such as unforcing it"
Clean-Text "suggested Magma presets button"
Clean-Text "whole measure- Pgup"
Clean-Text "whole measure- Pgdn"
Clean-Text "Starpower- u"
Clean-Text "BPM- i"
Clean-Text "GH3-style zones into Moonscraper please refer to"
Clean-Text "A. Moonscraper determines the length"
Clean-Text "If no audio is loaded Moonscraper shows a default"

# ------------------------------------------------------------------
# 2. The actual content fix: the "How do I place open notes?" answer.
#    Old key was 6; the manual now documents 0 (plus 1-5 also toggling
#    back to standard notes).
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "A. With the note tool selected press 6 on your keyboard. Press 6 on your keyboard again to return to standard notes. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "A. With the note tool selected press 0 on your keyboard. Press 0 on your keyboard again, or press any of the 1-5 keys to return to standard notes. ", `
    2) | Out-Null

# ------------------------------------------------------------------
# 3. Put the "_GoBack" bookmark back at the point of the edit above
#    (Word always parks it at the site of the most recent change).
# ------------------------------------------------------------------
$found = $d.Content
$found.Find.Execute("to return to standard notes", $false, $false, $false, `
                     $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($found.Start, $found.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ------------------------------------------------------------------
# 4. Minor proofing-language stamp Word adds to the flag-icon picture
#    run when it resaves the file (purely cosmetic, no visible text).
# ------------------------------------------------------------------
$picPara = $d.Paragraphs.Item(47)
$picRange = $picPara.Range
$picRange.LanguageID = "en-US"
$picRange.LanguageIDFarEast = "en-US"
